$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.840.20"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.092.69"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'233.84"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "'0.0784"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("E11").Value = "  +2.83%  "
$ws.Range("D12").Value = "'15.22"
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("D13").Value = "2.400.53"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "'21.38"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "2.085.79"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "37.794.68"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").Value = "'71.27"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "'230.54"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("D26").Value = "'9.83"
$ws.Range("E26").Value = "  +8.56%  "
$ws.Range("D27").Value = "'171.35"
$ws.Range("D28").Value = "'0.135"
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").Value = "'4.72"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("D35").Value = "'2.51"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").Value = "'3.32"
$ws.Range("E37").Value = "  -3.71%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("E40").Value = "  +8.93%  "
$ws.Range("D41").Value = "'102.05"
$ws.Range("E41").Value = "  +3.01%  "
$ws.Range("D42").Value = "'0.0975"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("E45").Value = "  +3.93%  "
$ws.Range("D46").Value = "1.448.83"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").Value = "'4.15"
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("D48").Value = "'1.07"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'7.27"
$ws.Range("E49").Value = "  -2.75%  "
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("D51").Value = "2.284.31"
$ws.Range("E51").Value = "  +0.12%  "
